$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (ratings) to be stored as text, matching the source data
# (ratings like "4.4" must not be auto-converted to numbers).
$ws.Range("D1:D24").NumberFormat = "@"

# Rows 17-24 are brand new (the sheet only went to row 16 before); give
# their date column (B) the same date/time display format used by the
# rest of column B, by copying the existing format instead of assigning
# a fresh NumberFormat (so it reuses the same underlying style as B1:B16).
$ws.Range("B1").Copy()
$ws.Range("B17:B24").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(1, 1).Value = "alle-bowls-acibadem-mah-kadikoy-istanbul"
$ws.Cells.Item(1, 2).Value = 44841.50505825232
$ws.Cells.Item(1, 3).Value = "AÇIK"
$ws.Cells.Item(1, 4).Value = "4.4"

$ws.Cells.Item(2, 1).Value = "ariana-s-cheesecake-acibadem-mah-kadikoy-istanbul"
$ws.Cells.Item(2, 2).Value = 44841.50512471065
$ws.Cells.Item(2, 3).Value = "KAPALI"
$ws.Cells.Item(2, 4).Value = "4.2"

$ws.Cells.Item(3, 1).Value = "restoran"
$ws.Cells.Item(3, 2).Value = 44841.50520675926
$ws.Cells.Item(3, 3).Value = "AÇIK"
$ws.Cells.Item(3, 4).Value = "4.3"

$ws.Cells.Item(4, 1).Value = "restoran"
$ws.Cells.Item(4, 2).Value = 44841.50527828703
$ws.Cells.Item(4, 3).Value = "AÇIK"
$ws.Cells.Item(4, 4).Value = "4.6"

$ws.Cells.Item(5, 1).Value = "alle-bowls-acibadem-mah-kadikoy-istanbul"
$ws.Cells.Item(5, 2).Value = 44841.50718738426
$ws.Cells.Item(5, 3).Value = "AÇIK"
$ws.Cells.Item(5, 4).Value = "4.4"

$ws.Cells.Item(6, 1).Value = "ariana-s-cheesecake-acibadem-mah-kadikoy-istanbul"
$ws.Cells.Item(6, 2).Value = 44841.50726765046
$ws.Cells.Item(6, 3).Value = "KAPALI"
$ws.Cells.Item(6, 4).Value = "4.2"

$ws.Cells.Item(7, 1).Value = "alle-bowls-acibadem-mah-kadikoy-istanbul"
$ws.Cells.Item(7, 2).Value = 44841.51303667824
$ws.Cells.Item(7, 3).Value = "AÇIK"
$ws.Cells.Item(7, 4).Value = "4.4"

$ws.Cells.Item(8, 1).Value = "ariana-s-cheesecake-acibadem-mah-kadikoy-istanbul"
$ws.Cells.Item(8, 2).Value = 44841.51311502315
$ws.Cells.Item(8, 3).Value = "KAPALI"
$ws.Cells.Item(8, 4).Value = "4.2"

$ws.Cells.Item(9, 1).Value = "restoran"
$ws.Cells.Item(9, 2).Value = 44841.51319752315
$ws.Cells.Item(9, 3).Value = "AÇIK"
$ws.Cells.Item(9, 4).Value = "4.3"

$ws.Cells.Item(10, 1).Value = "restoran"
$ws.Cells.Item(10, 2).Value = 44841.51327284722
$ws.Cells.Item(10, 3).Value = "AÇIK"
$ws.Cells.Item(10, 4).Value = "4.6"

$ws.Cells.Item(11, 1).Value = "cosa-bi-corba-bi-salata-acibadem-mah-kadikoy-istanbul"
$ws.Cells.Item(11, 2).Value = 44841.51334820602
$ws.Cells.Item(11, 3).Value = "AÇIK"
$ws.Cells.Item(11, 4).Value = "4.5"

$ws.Cells.Item(12, 1).Value = "restoran"
$ws.Cells.Item(12, 2).Value = 44841.51343730324
$ws.Cells.Item(12, 3).Value = "KAPALI"
$ws.Cells.Item(12, 4).Value = "4.4"

$ws.Cells.Item(13, 1).Value = "restoran"
$ws.Cells.Item(13, 2).Value = 44841.51351729166
$ws.Cells.Item(13, 3).Value = "KAPALI"
$ws.Cells.Item(13, 4).Value = "3.8"

$ws.Cells.Item(14, 1).Value = "doyuyo-sarayardi-cad-kadikoy-istanbul"
$ws.Cells.Item(14, 2).Value = 44841.51360233796
$ws.Cells.Item(14, 3).Value = "KAPALI"
$ws.Cells.Item(14, 4).Value = "4.5"

$ws.Cells.Item(15, 1).Value = "el-pollo-lasso-acibadem-mah-kadikoy-istanbul"
$ws.Cells.Item(15, 2).Value = 44841.51368598379
$ws.Cells.Item(15, 3).Value = "AÇIK"
$ws.Cells.Item(15, 4).Value = "4.3"

$ws.Cells.Item(16, 1).Value = "alle-bowls-acibadem-mah-kadikoy-istanbul"
$ws.Cells.Item(16, 2).Value = 44841.51379228009
$ws.Cells.Item(16, 3).Value = "AÇIK"
$ws.Cells.Item(16, 4).Value = "4.4"

$ws.Cells.Item(17, 1).Value = "ariana-s-cheesecake-acibadem-mah-kadikoy-istanbul"
$ws.Cells.Item(17, 2).Value = 44841.51386959491
$ws.Cells.Item(17, 3).Value = "KAPALI"
$ws.Cells.Item(17, 4).Value = "4.2"

$ws.Cells.Item(18, 1).Value = "alle-bowls-acibadem-mah-kadikoy-istanbul"
$ws.Cells.Item(18, 2).Value = 44841.51410226852
$ws.Cells.Item(18, 3).Value = "AÇIK"
$ws.Cells.Item(18, 4).Value = "4.4"

$ws.Cells.Item(19, 1).Value = "ariana-s-cheesecake-acibadem-mah-kadikoy-istanbul"
$ws.Cells.Item(19, 2).Value = 44841.51416706019
$ws.Cells.Item(19, 3).Value = "KAPALI"
$ws.Cells.Item(19, 4).Value = "4.2"

$ws.Cells.Item(20, 1).Value = "restoran"
$ws.Cells.Item(20, 2).Value = 44841.51423962963
$ws.Cells.Item(20, 3).Value = "AÇIK"
$ws.Cells.Item(20, 4).Value = "4.3"

$ws.Cells.Item(21, 1).Value = "restoran"
$ws.Cells.Item(21, 2).Value = 44841.51431548611
$ws.Cells.Item(21, 3).Value = "AÇIK"
$ws.Cells.Item(21, 4).Value = "4.6"

$ws.Cells.Item(22, 1).Value = "Alle Bowls"
$ws.Cells.Item(22, 2).Value = 44841.51753185411
$ws.Cells.Item(22, 3).Value = "AÇIK"
$ws.Cells.Item(22, 4).Value = "4.4"

$ws.Cells.Item(23, 1).Value = "Arianas Cheesecake"
$ws.Cells.Item(23, 2).Value = 44841.51761353836
$ws.Cells.Item(23, 3).Value = "KAPALI"
$ws.Cells.Item(23, 4).Value = "4.2"

$ws.Cells.Item(24, 1).Value = "Big Bold Quick"
$ws.Cells.Item(24, 2).Value = 44841.51770283478
$ws.Cells.Item(24, 3).Value = "AÇIK"
$ws.Cells.Item(24, 4).Value = "4.3"
